$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row reorder map: new row number -> old row number ---
$rowMap = @{
    2 = 3
    3 = 11
    4 = 2
    5 = 5
    6 = 4
    7 = 6
    8 = 13
    9 = 9
    10 = 7
    11 = 15
    12 = 8
    13 = 12
    14 = 10
    15 = 14
}

# --- Snapshot current A:I values for rows 2..15 before rewriting ---
$snapshot = @{}
for ($r = 2; $r -le 15; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 9; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# --- Write rows back out in the new order ---
for ($newR = 2; $newR -le 15; $newR++) {
    $oldR = $rowMap[$newR]
    $rowVals = $snapshot[$oldR]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($newR, $c).Value = $rowVals[$c - 1]
    }
}

# --- New header cells J1 / K1 ---
$ws.Cells.Item(1, 10).Value = 'Onkelos'
$ws.Cells.Item(1, 11).Value = 'Jonathan'

# --- New Targum columns (J = Onkelos, K = Jonathan) content ---
$jkData = @{
    4 = @('Anything that has a blemish on it you shall not bring, for it will not be accepted favorably for you.', 'But anything that hath a blemish you shall not offer; for that will not be acceptable from you.')
    6 = @('But if you will not do this, behold, you will have sinned against [<b>before</b>] Adonoy and you must realize that your sin will find you!', 'But if you will not perform this, behold, ye will have sinned before the Lord your God, and know that your sin will meet you.')
    8 = @('I see it but not now, I perceive it but not in the near future; a star [<b>king</b>] has gone forth from Yaakov, and a staff has arisen [<b>the Moshiach will be magnified by</b>] from Yisroel, which will smash the corners [<b>kill the leaders</b>] of Moav, and impale all of the sons of Sheis [<b>will rule over mankind</b>].', 'I shall see Him, but not now; I shall behold Him, but it is not near. When the mighty King of Jakob''s house shall reign, and the Meshiha, the Power-sceptre of Israel, be anointed, He will slay the princes of the Moabaee, and bring to nothing all the children of Sheth, the armies of Gog who will do battle against Israel and all their carcases shall fall before Him.')
    9 = @('He said, Please [<b>Now</b>] take your son, your only one, who you love—Yitzchok—and go to the land of Moriah [<b>worship</b>]. Sacrifice him [<b>before me</b>] as a burnt-offering on one of the mountains which I will designate to you.', 'And He said, Take now thy son, thy only one whom thou lovest, Izhak, and go into the land of worship, and offer him there, a whole burnt offering, upon one of the mountains that I will tell thee.')
    15 = @('A flawless lamb, a yearling male must be in your possession. You may take it from sheep or goats.', 'The lamb shall be perfect, a male, the son of a year he shall be to you; from the sheep or from the young goats ye may take.')
}
foreach ($r in $jkData.Keys) {
    $vals = $jkData[$r]
    $ws.Cells.Item($r, 10).Value = $vals[0]
    $ws.Cells.Item($r, 11).Value = $vals[1]
}

# --- Apply header style (bold/border, matches I1) to J1:K1 ---
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

# --- Apply body style (wrap text, matches I2:I15) to J2:K15 ---
$ws.Range("I2").Copy()
$ws.Range("J2:K15").PasteSpecial(-4122)

# --- Re-apply values after the format paste (PasteSpecial formats only, but be safe) ---
$ws.Cells.Item(1, 10).Value = 'Onkelos'
$ws.Cells.Item(1, 11).Value = 'Jonathan'
foreach ($r in $jkData.Keys) {
    $vals = $jkData[$r]
    $ws.Cells.Item($r, 10).Value = $vals[0]
    $ws.Cells.Item($r, 11).Value = $vals[1]
}

# --- Column widths for J (col 10) and K (col 11), matching target OOXML widths as closely as Excel rounding allows ---
$ws.Columns.Item(10).ColumnWidth = 391.5
$ws.Columns.Item(11).ColumnWidth = 440.83333333333337

$ws.Application.CutCopyMode = $false

